$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44435
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("S2").Value = 2000

# Row 4
$ws.Range("D4").Value = 44434
$ws.Range("M4").Value = 20

# Row 5
$ws.Range("D5").Value = 44511
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 28000
$ws.Range("S5").Value = 2800

# Row 7
$ws.Range("D7").Value = 44517
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 27000
$ws.Range("O7").Value = 27000
$ws.Range("P7").Value = 27000
$ws.Range("S7").Value = 2700

# Row 8
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 25000
$ws.Range("O8").Value = 25000
$ws.Range("P8").Value = 25000
$ws.Range("S8").Value = 2500

# Row 9
$ws.Range("D9").Value = 44476
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 2000

# Row 10
$ws.Range("D10").Value = 44466
$ws.Range("M10").Value = 60
